# Add 2022-Q4 data:
#  1) A new "2022-Q4" worksheet (fund-holdings detail), inserted right after
#     the "总计" summary sheet and before "2022-Q3".
#  2) A new row on the "总计" summary sheet referencing the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: create the "2022-Q4" worksheet.
#
# The easiest way to get an exact structural/style match with its sibling
# quarter sheets is to clone an existing one (it already carries the right
# header style, A-column style, column count, etc.) and then overwrite its
# data with the 2022-Q4 numbers.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $wb.Worksheets.Item("总计"))

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template ("2022-Q3") has 25 data rows; 2022-Q4 only needs 9, so drop
# the extra rows (11-26) entirely rather than just blanking them out.
$q4.Range("A11:H26").EntireRow.Delete()

# Fund holdings for 2022-Q4 (row, code, name, scale, position%, weight%, value, rank)
$q4Rows = @(
    @(2,  "014362", "睿远稳进配置两年持有混合A", "64.57", "37.12", "1.42", "0.9169", 7),
    @(3,  "014363", "睿远稳进配置两年持有混合C", "35.61", "37.12", "1.42", "0.5057", 7),
    @(4,  "001166", "建信环保产业股票",          "7.08",  "86.02", "3.74", "0.2648", 10),
    @(5,  "014522", "汇添富低碳投资一年持有混合A", "1.99",  "85.45", "2.29", "0.0456", 7),
    @(6,  "001484", "天弘新价值灵活配置混合A",    "1.26",  "94.05", "2.88", "0.0363", 6),
    @(7,  "002634", "华宝未来主导产业灵活配置混合A", "0.35",  "92.24", "8.04", "0.0281", 5),
    @(8,  "014523", "汇添富低碳投资一年持有混合C", "0.49",  "85.45", "2.29", "0.0112", 7),
    @(9,  "012919", "华宝未来主导产业灵活配置混合C", "0.07",  "92.24", "8.04", "0.0056", 5),
    @(10, "016246", "天弘新价值灵活配置混合C",    "0.03",  "94.05", "2.88", "0.0009", 6)
)

foreach ($r in $q4Rows) {
    $row   = $r[0]
    $code  = $r[1]
    $name  = $r[2]
    $scale = $r[3]
    $pos   = $r[4]
    $wgt   = $r[5]
    $val   = $r[6]
    $rank  = $r[7]

    # Fund code / numeric-looking figures are stored as text in this
    # workbook (leading zeros, fixed decimal display) - force text with a
    # leading apostrophe so Excel doesn't coerce them into numbers.
    $q4.Range("B$row").Value = "'" + $code
    $q4.Range("C$row").Value = $name
    $q4.Range("D$row").Value = "'" + $scale
    $q4.Range("E$row").Value = "'" + $pos
    $q4.Range("F$row").Value = "'" + $wgt
    $q4.Range("G$row").Value = "'" + $val
    $q4.Range("H$row").Value = $rank
}

# Copying a sheet makes the copy the active sheet (tabSelected); restore
# the selection to the real last sheet ("2020-Q4").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Select()

# ---------------------------------------------------------------------------
# Part 2: insert the 2022-Q4 summary row on "总计".
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Bring over the bordered/centered style used by the other index cells in
# column A (row 3 onward still carries it after the insert shifted them).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.82
